$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 18:22"

# 2) Update the per-country statistics that changed (values keyed to the
#    countries' CURRENT row positions before the re-sort below is applied).
#    Columns: B=Casos totales C=Nuevos casos D=Casos activos E=Recuperados
#             F=Casos criticos G=Muertes hoy H=Muertes

# Estados Unidos (row 4)
$ws.Range("B4").Value = 380698
$ws.Range("C4").Value = 13694
$ws.Range("D4").Value = 21309
$ws.Range("E4").Value = 347487
$ws.Range("F4").Value = 9104
$ws.Range("G4").Value = 1031
$ws.Range("H4").Value = 11902

# Italia (row 6)
$ws.Range("B6").Value = 135586
$ws.Range("C6").Value = 3039
$ws.Range("D6").Value = 24392
$ws.Range("E6").Value = 94067
$ws.Range("F6").Value = 3792
$ws.Range("G6").Value = 604
$ws.Range("H6").Value = 17127

# Alemania (row 7) - only E, G, H change
$ws.Range("E7").Value = 67536
$ws.Range("G7").Value = 92
$ws.Range("H7").Value = 1902

# Canada (row 16) - B, C, D, E change
$ws.Range("B16").Value = 17063
$ws.Range("C16").Value = 396
$ws.Range("D16").Value = 3796
$ws.Range("E16").Value = 12922

# Austria (row 17) - B, C, E change
$ws.Range("B17").Value = 12599
$ws.Range("C17").Value = 302
$ws.Range("E17").Value = 8310

# Chequia (row 30) - D, E, G, H change
$ws.Range("D30").Value = 172
$ws.Range("E30").Value = 4684
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 88

# Polonia (row 31) - B, C, E, G, H change
$ws.Range("B31").Value = 4666
$ws.Range("C31").Value = 253
$ws.Range("E31").Value = 4346
$ws.Range("G31").Value = 22
$ws.Range("H31").Value = 129

# Catar (row 49, before re-sort) - gets the new figures
$ws.Range("B49").Value = 2057
$ws.Range("C49").Value = 225
$ws.Range("D49").Value = 150
$ws.Range("E49").Value = 1901
$ws.Range("F49").Value = 37
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 6

# Argelia (row 57, before re-sort) - gets the new figures
$ws.Range("B57").Value = 1468
$ws.Range("C57").Value = 45
$ws.Range("D57").Value = 113
$ws.Range("E57").Value = 1162
$ws.Range("F57").Value = 46
$ws.Range("G57").Value = 20
$ws.Range("H57").Value = 193

# Republica de Chipre (row 86) - D, E, F change
$ws.Range("D86").Value = 47
$ws.Range("E86").Value = 438
$ws.Range("F86").Value = 13

# 3) Re-sort the whole country table (rows 4-215) by "Casos totales"
#    (column B) descending, which is how this sheet is always published -
#    the new figures shuffle country rank order.
$sortRange = $ws.Range("A4:H215")
$sortKey = $ws.Range("B4")
$sortRange.Sort($sortKey, 2)
